$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 and Row 3 (title/timestamp/historical distance/time bucket/uri) get
# swapped: the "Tornado History Project" entry now appears before the
# "Blizzard of '93" entry. Apply the new values directly (A:E) for rows 2 and 3.

$ws.Range("A2").Value = "Tornado History Project: March 13, 1993"
$ws.Range("B2").Value = "1993-03-13T00:00:00UTC"
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = "day_2_to_30"
$ws.Range("E2").Value = "http://www.tornadohistoryproject.com/tornado/1993/3/13/table"

$ws.Range("A3").Value = "Where were you during the Blizzard of '93? AL.com wants your pictures, memories"
$ws.Range("B3").Value = "2013-03-07T13:00:00UTC"
$ws.Range("C3").Value = 7311
$ws.Range("D3").Value = "day_31_beyond"
$ws.Range("E3").Value = "http://blog.al.com/spotnews/2013/03/where_were_you_during_the_bliz.html"
